$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 for the Fall_2018 semester (it becomes the new
# "current" semester). The previous row 2 (Spring_2018) shifts down to row 3
# and becomes "past", picking up the lastmod date it now needs.
$ws.Range("A2").EntireRow.Insert()

# New row 2: Fall_2018 / current / (no lastmod yet) / priority 0.6
$ws.Range("A2").Value = "Fall_2018"
$ws.Range("B2").Value = "current"
$ws.Range("D2").Value = 0.6

# Row 3 (previously row 2, Spring_2018) is now in the past with a lastmod
# date and the priority that used to belong to the most-recent past entry.
$ws.Range("B3").Value = "past"
$ws.Range("C3").Value = 43275
$ws.Range("D3").Value = 0.3

# Selection moved in the saved file.
$ws.Range("F7").Select()
